# Update on report progress.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Experiments (row 7): completion reached 100%, task note resolved/removed
$ws.Range("C7").Value = 1
$ws.Range("D7").ClearContents()

# Bibliography (row 12): completion reached 100%, task note resolved/removed
$ws.Range("C12").Value = 1
$ws.Range("D12").ClearContents()

# Appendix (row 13): completion reached 100%, task note resolved/removed
$ws.Range("C13").Value = 1
$ws.Range("D13").ClearContents()

# Recalculate dependent formulas (B14 total, C14 average completion)
$excel.Calculate()

# Update the last selected cell/cursor position on the sheet
$ws.Range("D18").Select()
